$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT (t="inlineStr"/shared-string), even when
# the text looks like a number ("1", "2", "76", ...), without leaving a
# permanent numeric-format style on the cell (format is applied only long
# enough to force the text type, then cleared again).
function Set-TextCell($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Header J1: "# Jaula" -> "numeroJaula"
$ws.Range("J1").Value = "numeroJaula"

# Full data set for rows 2-9 (ID, Cuerda, Frente, Anillo, Placa, Color, Peso, Ciudad, Tipo, numeroJaula)
$data = @(
    @{ id=1; cuerda="Alacranes"; frente="2";   anillo="1";   placa="2";   color="giro";     peso=3;   ciudad="sogamoso"; tipo="Gallo"; jaula="1" }
    @{ id=2; cuerda="Alacranes"; frente="2";   anillo="2";   placa="3";   color="colorado";  peso=3.2; ciudad="sogamoso"; tipo="Pollo"; jaula="2" }
    @{ id=3; cuerda="Sara";      frente="1";   anillo="3";   placa="4";   color="cenizo";    peso=3.1; ciudad="tunja";    tipo="Gallo"; jaula="3" }
    @{ id=4; cuerda="Sara";      frente="1";   anillo="5";   placa="6";   color="javado";    peso=3.3; ciudad="tunja";    tipo="Pollo"; jaula="4" }
    @{ id=5; cuerda="Alacranes"; frente="1";   anillo="76";  placa="98";  color="colorado";  peso=3.4; ciudad="bogota";   tipo="Gallo"; jaula="5" }
    @{ id=6; cuerda="Alacranes"; frente="1";   anillo="88";  placa="99";  color="canaguai";  peso=3.8; ciudad="bogota";   tipo="Gallo"; jaula="6" }
    @{ id=7; cuerda="Nathaly";   frente="1";   anillo="826"; placa="563"; color="giro";      peso=3.6; ciudad="Duitama";  tipo="Pollo"; jaula="7" }
    @{ id=8; cuerda="Nathaly";   frente="1";   anillo="21";  placa="34";  color="giro";      peso=3.9; ciudad="Duitama";  tipo="Gallo"; jaula="8" }
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row.id          # A: ID (number)
    $ws.Range("B$r").Value = $row.cuerda            # B: Cuerda (text)
    Set-TextCell $ws.Range("C$r") $row.frente       # C: Frente (numeric-looking text)
    Set-TextCell $ws.Range("D$r") $row.anillo       # D: Anillo (numeric-looking text)
    Set-TextCell $ws.Range("E$r") $row.placa        # E: Placa (numeric-looking text)
    $ws.Range("F$r").Value = $row.color             # F: Color (text)
    $ws.Cells.Item($r, 7).Value = $row.peso          # G: Peso (number)
    $ws.Range("H$r").Value = $row.ciudad            # H: Ciudad (text)
    $ws.Range("I$r").Value = $row.tipo              # I: Tipo (text)
    Set-TextCell $ws.Range("J$r") $row.jaula        # J: numeroJaula (numeric-looking text)
    $r++
}
